$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = [double]"0.6096806208067198"
$ws.Range("I2").Value = [double]"0.6096806208067198"
$ws.Range("L2").Value = [double]"2.873323447694827"
$ws.Range("M2").Value = '[-5.235292234706125, 10.981939130095778]'
$ws.Range("N2").Value = [double]"0.4790942653775361"
$ws.Range("O2").Value = [double]"0.4790942653775361"
$ws.Range("P2").Value = [double]"-0.8176317216550011"
$ws.Range("Q2").Value = '[-3.956079637853814, 2.320816194543812]'
$ws.Range("R2").Value = [double]"0.6023559687572058"
$ws.Range("S2").Value = [double]"0.6023559687572058"
$ws.Range("T2").Value = [double]"11.27073123817208"
$ws.Range("U2").Value = '[7.059477421560091, 15.481985054784067]'
$ws.Range("V2").Value = [double]"2.482267692949236e-06"
$ws.Range("W2").Value = [double]"2.482267692949236e-06"
$ws.Range("X2").Value = [double]"3.382082082082164"
$ws.Range("Y2").Value = [double]"-9.599909909910142"
$ws.Range("Z2").Value = [double]"16.36407407407447"

# Row 3
$ws.Range("H3").Value = [double]"0.4364407198044931"
$ws.Range("I3").Value = [double]"0.4364407198044931"
$ws.Range("L3").Value = [double]"4.189797053738692"
$ws.Range("M3").Value = '[-4.333231211233798, 12.712825318711182]'
$ws.Range("N3").Value = [double]"0.3274174498818065"
$ws.Range("O3").Value = [double]"0.3274174498818065"
$ws.Range("P3").Value = [double]"-1.396263401595464"
$ws.Range("Q3").Value = '[-4.50955341866643, 1.7170266154755023]'
$ws.Range("R3").Value = [double]"0.3711796174488051"
$ws.Range("S3").Value = [double]"0.3711796174488051"
$ws.Range("T3").Value = [double]"11.87289968298832"
$ws.Range("U3").Value = '[7.073513567262134, 16.672285798714512]'
$ws.Range("V3").Value = [double]"9.748779093099458e-06"
$ws.Range("W3").Value = [double]"9.748779093099458e-06"
$ws.Range("X3").Value = [double]"5.775555555555695"
$ws.Range("Y3").Value = [double]"-7.102372372372542"
$ws.Range("Z3").Value = [double]"18.65348348348393"

# Row 4
$ws.Range("H4").Value = [double]"0.7416111861138102"
$ws.Range("I4").Value = [double]"0.7416111861138102"
$ws.Range("L4").Value = [double]"3.160552264796364"
$ws.Range("M4").Value = '[-7.847627953370742, 14.168732482963469]'
$ws.Range("N4").Value = [double]"0.5659660593504228"
$ws.Range("O4").Value = [double]"0.5659660593504228"
$ws.Range("P4").Value = [double]"-1.836526636332772"
$ws.Range("Q4").Value = '[-4.9309482290578535, 1.2578949563923087]'
$ws.Range("R4").Value = [double]"0.2382071584112135"
$ws.Range("S4").Value = [double]"0.2382071584112135"
$ws.Range("T4").Value = [double]"15.21563789837347"
$ws.Range("U4").Value = '[9.29253624748165, 21.138739549265285]'
$ws.Range("V4").Value = [double]"5.143654590478164e-06"
$ws.Range("W4").Value = [double]"5.143654590478164e-06"
$ws.Range("X4").Value = [double]"7.596676676676861"
$ws.Range("Y4").Value = [double]"-5.203203203203325"
$ws.Range("Z4").Value = [double]"20.39655655655704"

# Row 5
$ws.Range("B5").Value = [double]"1"
$ws.Range("H5").Value = [double]"0.007290346212915022"
$ws.Range("I5").Value = [double]"0.007290346212915022"
$ws.Range("L5").Value = [double]"9.860930522735769"
$ws.Range("M5").Value = '[1.8026379472797913, 17.919223098191747]'
$ws.Range("N5").Value = [double]"0.01759757688795527"
$ws.Range("O5").Value = [double]"0.01759757688795527"
$ws.Range("P5").Value = [double]"-2.025210879791619"
$ws.Range("Q5").Value = '[-2.9183162988301588, -1.1321054607530785]'
$ws.Range("R5").Value = [double]"3.827129155076214e-05"
$ws.Range("S5").Value = [double]"3.827129155076214e-05"
$ws.Range("T5").Value = [double]"13.88517860085352"
$ws.Range("U5").Value = '[9.625373566210166, 18.144983635496885]'
$ws.Range("V5").Value = [double]"4.510260587231585e-08"
$ws.Range("W5").Value = [double]"4.510260587231585e-08"
$ws.Range("X5").Value = [double]"8.37715715715736"
$ws.Range("Y5").Value = [double]"4.682882882882995"
$ws.Range("Z5").Value = [double]"12.07143143143172"

# Row 6
$ws.Range("F6").Value = [double]"22.39000000000006"
$ws.Range("H6").Value = [double]"0.03560842243410001"
$ws.Range("I6").Value = [double]"0.03560842243410001"
$ws.Range("L6").Value = [double]"8.19964008450466"
$ws.Range("M6").Value = '[-0.22497872915635853, 16.624258898165678]'
$ws.Range("N6").Value = [double]"0.05616766529417516"
$ws.Range("O6").Value = [double]"0.05616766529417516"
$ws.Range("P6").Value = [double]"2.584974135386196"
$ws.Range("Q6").Value = '[1.0377633390236554, 4.132184931748737]'
$ws.Range("R6").Value = [double]"0.001573373014130297"
$ws.Range("S6").Value = [double]"0.001573373014130297"
$ws.Range("T6").Value = [double]"13.19005232142028"
$ws.Range("U6").Value = '[8.740115940005897, 17.639988702834664]'
$ws.Range("V6").Value = [double]"3.459460213850463e-07"
$ws.Range("W6").Value = [double]"3.459460213850463e-07"
$ws.Range("X6").Value = [double]"13.17849849849853"
$ws.Range("Y6").Value = [double]"7.665045045045063"
$ws.Range("Z6").Value = [double]"18.691951951952"

# Row 7
$ws.Range("B7").Value = [double]"0"
$ws.Range("F7").Value = [double]"22.39000000000006"
$ws.Range("H7").Value = [double]"0.06925100668920903"
$ws.Range("I7").Value = [double]"0.06925100668920903"
$ws.Range("L7").Value = [double]"7.495142413179052"
$ws.Range("M7").Value = '[-1.2064155497007834, 16.19670037605889]'
$ws.Range("N7").Value = [double]"0.08961079478189071"
$ws.Range("O7").Value = [double]"0.08961079478189071"
$ws.Range("P7").Value = [double]"2.647868883205812"
$ws.Range("Q7").Value = '[0.7736053981812714, 4.522132368230352]'
$ws.Range("R7").Value = [double]"0.006654147679102662"
$ws.Range("S7").Value = [double]"0.006654147679102662"
$ws.Range("T7").Value = [double]"14.3046326107477"
$ws.Range("U7").Value = '[9.711391600364117, 18.897873621131275]'
$ws.Range("V7").Value = [double]"1.229012200898438e-07"
$ws.Range("W7").Value = [double]"1.229012200898438e-07"
$ws.Range("X7").Value = [double]"12.95437437437441"
$ws.Range("Y7").Value = [double]"6.275475475475494"
$ws.Range("Z7").Value = [double]"19.63327327327332"

# Row 8
$ws.Range("F8").Value = [double]"22.39000000000006"
$ws.Range("H8").Value = [double]"0.1838156260647319"
$ws.Range("I8").Value = [double]"0.1838156260647319"
$ws.Range("L8").Value = [double]"6.254283454279176"
$ws.Range("M8").Value = '[-2.8134449208877887, 15.32201182944614]'
$ws.Range("N8").Value = [double]"0.1716122905644439"
$ws.Range("O8").Value = [double]"0.1716122905644439"
$ws.Range("P8").Value = [double]"-2.44031621540108"
$ws.Range("Q8").Value = '[-5.566185182035969, 0.6855527512338089]'
$ws.Range("R8").Value = [double]"0.1228682772488165"
$ws.Range("S8").Value = [double]"0.1228682772488165"
$ws.Range("T8").Value = [double]"14.0633226348365"
$ws.Range("U8").Value = '[9.326935052121737, 18.799710217551254]'
$ws.Range("V8").Value = [double]"3.339920906508809e-07"
$ws.Range("W8").Value = [double]"3.339920906508809e-07"
$ws.Range("X8").Value = [double]"8.696016016016038"
$ws.Range("Y8").Value = [double]"-2.442952952952961"
$ws.Range("Z8").Value = [double]"19.83498498498504"

# Row 9
$ws.Range("F9").Value = [double]"22.39000000000006"
$ws.Range("H9").Value = [double]"0.3072117336260169"
$ws.Range("I9").Value = [double]"0.3072117336260169"
$ws.Range("L9").Value = [double]"5.477257122920014"
$ws.Range("M9").Value = '[-4.006836568663264, 14.961350814503291]'
$ws.Range("N9").Value = [double]"0.2508852380095776"
$ws.Range("O9").Value = [double]"0.2508852380095776"
$ws.Range("P9").Value = [double]"-2.553526761476388"
$ws.Range("Q9").Value = '[-5.6856852028932385, 0.5786316799404627]'
$ws.Range("R9").Value = [double]"0.107557078832887"
$ws.Range("S9").Value = [double]"0.107557078832887"
$ws.Range("T9").Value = [double]"14.62346983513498"
$ws.Range("U9").Value = '[9.592406595076064, 19.654533075193893]'
$ws.Range("V9").Value = [double]"5.136056155929936e-07"
$ws.Range("W9").Value = [double]"5.136056155929936e-07"
$ws.Range("X9").Value = [double]"9.099439439439463"
$ws.Range("Y9").Value = [double]"-2.061941941941948"
$ws.Range("Z9").Value = [double]"20.26082082082088"

# Row 10
$ws.Range("F10").Value = [double]"22.39000000000006"
$ws.Range("H10").Value = [double]"0.3534952714344888"
$ws.Range("I10").Value = [double]"0.3534952714344888"
$ws.Range("L10").Value = [double]"5.859458799627342"
$ws.Range("M10").Value = '[-4.905805663025955, 16.62472326228064]'
$ws.Range("N10").Value = [double]"0.2787961898415967"
$ws.Range("O10").Value = [double]"0.2787961898415967"
$ws.Range("P10").Value = [double]"-2.025210879791619"
$ws.Range("Q10").Value = '[-5.157369321208469, 1.1069475616252316]'
$ws.Range("R10").Value = [double]"0.1994429032139815"
$ws.Range("S10").Value = [double]"0.1994429032139815"
$ws.Range("T10").Value = [double]"15.73925255660278"
$ws.Range("U10").Value = '[10.160560820431481, 21.317944292774072]'
$ws.Range("V10").Value = [double]"9.223455972318106e-07"
$ws.Range("W10").Value = [double]"9.223455972318106e-07"
$ws.Range("X10").Value = [double]"7.216796796796819"
$ws.Range("Y10").Value = [double]"-3.944584584584591"
$ws.Range("Z10").Value = [double]"18.37817817817823"
